# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 7
    3  = 8
    4  = 7
    5  = 4
    6  = 7
    7  = 8
    8  = 6
    9  = 9
    10 = 6
    11 = 8
    12 = 7
    13 = 10
    14 = 10
    15 = 7
    16 = 10
    17 = 3
    18 = 6
    19 = 5
    20 = 8
    21 = 11
    22 = 2
    23 = 7
    24 = 11
    25 = 6
    26 = 4
    27 = 6
    28 = 4
    29 = 5
    30 = 4
    31 = 7
    32 = 5
    33 = 4
    34 = 2
    35 = 4
    36 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
